{"js": "// Word Online / Office.js edit script.\n// Applies the \"store confirm Schedule update Skype meetings\" change:\n//  1. \"Mikaeil\" (in the Attended-by line) gets struck through, splitting the\n//     run into a plain leading space and a struck-through name.\n//  2. The \"Planned activities:\" paragraph loses the _GoBack bookmark (it is\n//     relocated further down).\n//  3. Four new paragraphs of meeting notes are inserted right after\n//     \"Planned activities:\", followed by the pre-existing \"Highlights:\"\n//     paragraph content (re-created fresh) and the old trailing\n//     \"Highlights:\" paragraph is removed once its text has been\n//     reproduced by the freshly inserted one.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load('text');\nawait context.sync();\n\n// ---- 1. Strike through \"Mikaeil\" -----------------------------------------\nconst mikaeil = body.search(\"Mikaeil\", { matchCase: true });\nmikaeil.load('text');\nawait context.sync();\nmikaeil.items[0].font.strikeThrough = true;\nawait context.sync();\n\n// ---- 2. Locate the anchor paragraphs --------------------------------------\nparagraphs.load('text');\nawait context.sync();\n\nlet plannedPara = null;\nlet oldHighlightsPara = null;\nfor (const p of paragraphs.items) {\n  const t = p.text.trim();\n  if (t === \"Planned activities:\" && !plannedPara) {\n    plannedPara = p;\n  } else if (t === \"Highlights:\") {\n    oldHighlightsPara = p;\n  }\n}\n\n// The bookmark currently sits at the end of \"Planned activities:\" \u2014 it moves\n// further down in the new text, so drop it from here first.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// ---- 3. Insert the new paragraphs after \"Planned activities:\" -------------\nconst pHighlights = plannedPara.insertParagraph(\"Highlights:\", \"After\");\nconst pProgress = pHighlights.insertParagraph(\"How is progress\", \"After\");\n\n// \"Website functionality : Tickets, Booking\" + line break + \"Process of SQL\"\nconst pWebsite = pProgress.insertParagraph(\"\", \"After\");\nconst websiteOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>\n              <w:t>Website functionality : Tickets, Booking</w:t>\n            </w:r>\n            <w:r>\n              <w:br/>\n              <w:t>Process of SQL</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\npWebsite.insertOoxml(websiteOoxml, \"Replace\");\n\n// \"Discussion on additional division of website and application\" + bookmark + \" tasks\"\nconst pDiscussion = pWebsite.insertParagraph(\"\", \"After\");\nconst discussionOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>\n              <w:t xml:space=\"preserve\">Discussion on additional division of website and </w:t>\n            </w:r>\n            <w:r>\n              <w:t>application</w:t>\n            </w:r>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n            <w:r>\n              <w:t xml:space=\"preserve\"> tasks</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\npDiscussion.insertOoxml(discussionOoxml, \"Replace\");\n\n// ---- 4. Drop the original trailing \"Highlights:\" paragraph ----------------\n// (its text now lives in the freshly inserted pHighlights paragraph above).\nif (oldHighlightsPara) {\n  oldHighlightsPara.delete();\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# Applies the \"store confirm Schedule update Skype meetings\" change:\n#  1. \"Mikaeil\" (in the Attended-by line) gets struck through, splitting the\n#     run into a plain leading space and a struck-through name.\n#  2. The \"Planned activities:\" paragraph loses the _GoBack bookmark (it is\n#     relocated further down, to the end of the new notes).\n#  3. Four new paragraphs of meeting notes are inserted right after\n#     \"Planned activities:\" (Highlights:, How is progress, the\n#     Website-functionality line with an embedded line break, and the\n#     Discussion line carrying the relocated bookmark).\n#  4. The original trailing \"Highlights:\" paragraph \u2014 whose text now lives\n#     in the freshly inserted paragraph from step 3 \u2014 is removed.\n\n$d = $word.ActiveDocument\n\n# ---- 1. Strike through \"Mikaeil\" ------------------------------------------\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"Mikaeil\"\n$find.Execute() | Out-Null\n$rng.Font.StrikeThrough = 1\n\n# ---- 2. Remove the _GoBack bookmark from its current location -------------\n$d.Bookmarks(\"_GoBack\").Delete()\n\n# ---- 3. Locate the \"Planned activities:\" paragraph and insert after it ----\n$plannedIndex = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text.Trim() -eq \"Planned activities:\") {\n        $plannedIndex = $i\n        break\n    }\n}\n\n$p = $d.Paragraphs($plannedIndex)\n$p.Range.InsertParagraphAfter()\n$p = $d.Paragraphs($plannedIndex + 1)\n$p.Range.Text = \"Highlights:\"\n\n$p.Range.InsertParagraphAfter()\n$p = $d.Paragraphs($plannedIndex + 2)\n$p.Range.Text = \"How is progress\"\n\n$p.Range.InsertParagraphAfter()\n$p = $d.Paragraphs($plannedIndex + 3)\n$websiteXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:t>Website functionality : Tickets, Booking</w:t></w:r><w:r><w:br/><w:t>Process of SQL</w:t></w:r></w:p>'\n$p.Range.InsertXML($websiteXml) | Out-Null\n\n$p.Range.InsertParagraphAfter()\n$p = $d.Paragraphs($plannedIndex + 4)\n$discussionXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:t xml:space=\"preserve\">Discussion on additional division of website and </w:t></w:r><w:r><w:t>application</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/><w:r><w:t xml:space=\"preserve\"> tasks</w:t></w:r></w:p>'\n$p.Range.InsertXML($discussionXml) | Out-Null\n\n# ---- 4. Drop the original trailing \"Highlights:\" paragraph ----------------\n# (its text now lives in the freshly inserted paragraph from step 3 above).\n$last = $d.Paragraphs($d.Paragraphs.Count)\nif ($last.Range.Text.Trim() -eq \"Highlights:\") {\n    $last.Range.Delete()\n}\n"}
